$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was collected for "Arveja Verde" (Femacal de
# La Calera). It belongs chronologically before the existing row 21, so
# insert a new row at position 21 and push the rest of the table down.
$ws.Rows.Item(21).Insert()

$ws.Range("A21").Value = 3
$ws.Range("B21").Value = "Femacal de La Calera"
$ws.Range("C21").Value = "Coquimbo"
$ws.Range("D21").Value = 44510
$ws.Range("D21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 100112022
$ws.Range("G21").Value = "Arveja Verde"
$ws.Range("H21").Value = "Perfection"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 73
$ws.Range("K21").Value = 16500
$ws.Range("L21").Value = 17000
$ws.Range("M21").Value = 16740
$ws.Range("N21").Value = "$/saco 25 kilos"
$ws.Range("O21").Value = "Provincia de Limarí"
$ws.Range("P21").Value = 670
$ws.Range("Q21").Value = 25
$ws.Range("R21").Value = "Hortaliza"
